$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

$newCC = "Julian Andres Pena Ospina <japenao@incauca.com>; Ginna Constanza Rosero Arevalo <gcrosero@incauca.com>; Félix Andrés  Molina Serrano <famolina@incauca.com>; Aderson Orozco Gonzalez <aorozco@incauca.com>; Edwin Fabian Mesias <efmesias@incauca.com>"

for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $newCC
}

$ws.Range("J5").Font.Italic = $true

$ws.Range("D1").Select()
